# Update "想去人数" (F column) counts for several anime-expo rows, as
# regenerated by the gh-pages data refresh (commit 456a3b4).
#
# The workbook contains the same events duplicated on the "展览" sheet and
# the "全部类型" sheet (which additionally includes one extra 演出 row, so
# row numbers from row 14 onward are shifted by +1 relative to "展览").

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{
        3  = 1355
        7  = 11674
        8  = 4385
        10 = 38
        15 = 146
        17 = 5093
        18 = 63
        20 = 514
        21 = 11336
        22 = 11266
        27 = 47
    }
    "全部类型" = @{
        3  = 1355
        7  = 11674
        8  = 4385
        10 = 38
        16 = 146
        18 = 5093
        19 = 63
        21 = 514
        22 = 11336
        23 = 11266
        28 = 47
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowUpdates = $updates[$sheetName]
    foreach ($row in $rowUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowUpdates[$row]
    }
}
